# Auto-generated edits applying the diff to Kujata_Profits workbook
# (market-data refresh across the per-job sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2299.6155
$ws.Range("I51").Value = 1965
$ws.Range("J51").Value = 2400
$ws.Range("K51").Value = 1965
$ws.Range("L51").Value = 2400
$ws.Range("M51").Value = -1481
$ws.Range("N51").Value = -3368
$ws.Range("H113").Value = 3616
$ws.Range("I113").Value = 3578.25
$ws.Range("K113").Value = 3578.25
$ws.Range("M113").Value = -324.25
$ws.Range("H116").Value = 2775
$ws.Range("I116").Value = 1934
$ws.Range("K116").Value = 1934
$ws.Range("M116").Value = 1508
$ws.Range("H132").Value = 10424794
$ws.Range("I132").Value = 14500302
$ws.Range("J132").Value = 9606.888999999999
$ws.Range("K132").Value = 43500906
$ws.Range("L132").Value = 28820.667
$ws.Range("M132").Value = -43498376
$ws.Range("N132").Value = -33880.667
$ws.Range("H137").Value = 1547.7667
$ws.Range("I137").Value = 961.73334
$ws.Range("J137").Value = 2133.8
$ws.Range("K137").Value = 2885.20002
$ws.Range("L137").Value = 6401.400000000001
$ws.Range("M137").Value = -335.2000200000002
$ws.Range("N137").Value = -11501.4
$ws.Range("H138").Value = 1145.1
$ws.Range("I138").Value = 618.6087
$ws.Range("J138").Value = 1593.5927
$ws.Range("K138").Value = 1855.8261
$ws.Range("L138").Value = 4780.7781
$ws.Range("M138").Value = 3284.1739
$ws.Range("N138").Value = -15060.7781
$ws.Range("H141").Value = 616.7273
$ws.Range("I141").Value = 616.7273
$ws.Range("K141").Value = 1850.1819
$ws.Range("M141").Value = 3329.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4077.4026
$ws.Range("I32").Value = 3512.9297
$ws.Range("J32").Value = 10757
$ws.Range("K32").Value = 3512.9297
$ws.Range("L32").Value = 10757
$ws.Range("M32").Value = -3225.9297
$ws.Range("N32").Value = -11331
$ws.Range("H61").Value = 25641878
$ws.Range("I61").Value = 30303734
$ws.Range("J61").Value = 1669
$ws.Range("K61").Value = 30303734
$ws.Range("L61").Value = 1669
$ws.Range("M61").Value = -30303522
$ws.Range("N61").Value = -2093
$ws.Range("H63").Value = 28573474
$ws.Range("I63").Value = 2009.7391
$ws.Range("J63").Value = 83335450
$ws.Range("K63").Value = 2009.7391
$ws.Range("L63").Value = 83335450
$ws.Range("M63").Value = -1323.7391
$ws.Range("N63").Value = -83336822
$ws.Range("H66").Value = 28573474
$ws.Range("I66").Value = 2009.7391
$ws.Range("J66").Value = 83335450
$ws.Range("K66").Value = 10048.6955
$ws.Range("L66").Value = 416677250
$ws.Range("M66").Value = -6616.6955
$ws.Range("N66").Value = -416684114
$ws.Range("H74").Value = 2230.1667
$ws.Range("I74").Value = 1593.5
$ws.Range("K74").Value = 1593.5
$ws.Range("M74").Value = -719.5
$ws.Range("H77").Value = 2230.1667
$ws.Range("I77").Value = 1593.5
$ws.Range("K77").Value = 7967.5
$ws.Range("M77").Value = -3599.5
$ws.Range("H97").Value = 256.38235
$ws.Range("I97").Value = 281.4643
$ws.Range("J97").Value = 139.33333
$ws.Range("K97").Value = 281.4643
$ws.Range("L97").Value = 139.33333
$ws.Range("M97").Value = 214.5357
$ws.Range("N97").Value = -1131.33333
$ws.Range("H132").Value = 1027.5
$ws.Range("I132").Value = 948.9091
$ws.Range("J132").Value = 1243.625
$ws.Range("K132").Value = 2846.7273
$ws.Range("L132").Value = 3730.875
$ws.Range("M132").Value = -316.7273
$ws.Range("N132").Value = -8790.875
$ws.Range("H136").Value = 25641878
$ws.Range("I136").Value = 30303734
$ws.Range("J136").Value = 1669
$ws.Range("K136").Value = 90911202
$ws.Range("L136").Value = 5007
$ws.Range("M136").Value = -90908652
$ws.Range("N136").Value = -10107

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 62501132
$ws.Range("I99").Value = 83334184
$ws.Range("K99").Value = 83334184
$ws.Range("M99").Value = -83332686
$ws.Range("H105").Value = 58825452
$ws.Range("I105").Value = 71430300
$ws.Range("J105").Value = 2837
$ws.Range("K105").Value = 71430300
$ws.Range("L105").Value = 2837
$ws.Range("M105").Value = -71428553
$ws.Range("N105").Value = -6331
$ws.Range("H134").Value = 3237.0208
$ws.Range("I134").Value = 940.38635
$ws.Range("J134").Value = 28500
$ws.Range("K134").Value = 2821.15905
$ws.Range("L134").Value = 85500
$ws.Range("M134").Value = -286.1590500000002
$ws.Range("N134").Value = -90570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 55556868
$ws.Range("I16").Value = 76924340
$ws.Range("J16").Value = 1465
$ws.Range("K16").Value = 76924340
$ws.Range("L16").Value = 1465
$ws.Range("M16").Value = -76924053
$ws.Range("N16").Value = -2039
$ws.Range("H31").Value = 2126.25
$ws.Range("I31").Value = 1945.8462
$ws.Range("J31").Value = 2908
$ws.Range("K31").Value = 1945.8462
$ws.Range("L31").Value = 2908
$ws.Range("M31").Value = -1650.8462
$ws.Range("N31").Value = -3498
$ws.Range("H34").Value = 2126.25
$ws.Range("I34").Value = 1945.8462
$ws.Range("J34").Value = 2908
$ws.Range("K34").Value = 1945.8462
$ws.Range("L34").Value = 2908
$ws.Range("M34").Value = -1743.8462
$ws.Range("N34").Value = -3312
$ws.Range("H58").Value = 672.48334
$ws.Range("I58").Value = 580.71155
$ws.Range("K58").Value = 580.71155
$ws.Range("M58").Value = -377.71155
$ws.Range("H107").Value = 671.0625
$ws.Range("J107").Value = 898.3333
$ws.Range("L107").Value = 898.3333
$ws.Range("N107").Value = -4738.3333
$ws.Range("H113").Value = 55556868
$ws.Range("I113").Value = 76924340
$ws.Range("J113").Value = 1465
$ws.Range("K113").Value = 76924340
$ws.Range("L113").Value = 1465
$ws.Range("M113").Value = -76922170
$ws.Range("N113").Value = -5805
$ws.Range("H132").Value = 3053.276
$ws.Range("I132").Value = 3325.24
$ws.Range("J132").Value = 1353.5
$ws.Range("K132").Value = 9975.719999999999
$ws.Range("L132").Value = 4060.5
$ws.Range("M132").Value = -7445.719999999999
$ws.Range("N132").Value = -9120.5
$ws.Range("H134").Value = 9092072
$ws.Range("I134").Value = 1199.3617
$ws.Range("K134").Value = 3598.0851
$ws.Range("M134").Value = -1063.0851
$ws.Range("H136").Value = 672.48334
$ws.Range("I136").Value = 580.71155
$ws.Range("K136").Value = 1742.13465
$ws.Range("M136").Value = 807.86535
$ws.Range("H138").Value = 127796.664
$ws.Range("J138").Value = 127796.664
$ws.Range("L138").Value = 127796.664
$ws.Range("N138").Value = -138076.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 617.6977000000001
$ws.Range("J113").Value = 647.7222
$ws.Range("L113").Value = 1943.1666
$ws.Range("N113").Value = -6283.1666
$ws.Range("H121").Value = 603.2
$ws.Range("J121").Value = 722
$ws.Range("L121").Value = 2166
$ws.Range("N121").Value = -4786
$ws.Range("H131").Value = 27779212
$ws.Range("J131").Value = 1679.5358
$ws.Range("L131").Value = 5038.607400000001
$ws.Range("N131").Value = -15118.6074

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1687.1666
$ws.Range("H127").Value = 37133.332
$ws.Range("J127").Value = 37133.332
$ws.Range("L127").Value = 37133.332
$ws.Range("N127").Value = -47053.332
$ws.Range("H132").Value = 1838.7407
$ws.Range("I132").Value = 1529.3636
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 4588.0908
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -2058.0908
$ws.Range("N132").Value = -14660

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2713.8
$ws.Range("I40").Value = 2459.7778
$ws.Range("K40").Value = 2459.7778
$ws.Range("M40").Value = -2323.7778
$ws.Range("H132").Value = 29307.111
$ws.Range("I132").Value = 1438.3667
$ws.Range("J132").Value = 168650.83
$ws.Range("K132").Value = 4315.1001
$ws.Range("L132").Value = 505952.49
$ws.Range("M132").Value = -1785.1001
$ws.Range("N132").Value = -511012.49
$ws.Range("H136").Value = 5382.84
$ws.Range("I136").Value = 5482.125
$ws.Range("K136").Value = 16446.375
$ws.Range("M136").Value = -13896.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 583.4
$ws.Range("I107").Value = 504.85715
$ws.Range("J107").Value = 766.6667
$ws.Range("K107").Value = 1514.57145
$ws.Range("L107").Value = 2300.0001
$ws.Range("M107").Value = 405.4285500000001
$ws.Range("N107").Value = -6140.0001
$ws.Range("H113").Value = 305.91666
$ws.Range("I113").Value = 201.3
$ws.Range("J113").Value = 380.64285
$ws.Range("K113").Value = 603.9000000000001
$ws.Range("L113").Value = 1141.92855
$ws.Range("M113").Value = 1566.1
$ws.Range("N113").Value = -5481.928550000001
$ws.Range("H132").Value = 1828.3478
$ws.Range("I132").Value = 2328.3635
$ws.Range("J132").Value = 1370
$ws.Range("K132").Value = 6985.0905
$ws.Range("L132").Value = 4110
$ws.Range("M132").Value = -4455.0905
$ws.Range("N132").Value = -9170
$ws.Range("H136").Value = 539.2
$ws.Range("I136").Value = 539.2
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 1617.6
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 932.3999999999999
$ws.Range("N136").ClearContents()
$ws.Range("H138").Value = 34793.332
$ws.Range("J138").Value = 34793.332
$ws.Range("L138").Value = 34793.332
$ws.Range("N138").Value = -45073.332

